# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first worker row (ADRIANA LUCIA LOPEZ ACEVEDO / 1047501495) entirely;
# the remaining workers (YESICA, DALYS) shift up one row.
$ws.Rows.Item(16).Delete()

# Update the aggregate "VALOR MORA" total.
$ws.Range("E11").Value = 456940

# Update the worker count (Cant. Trabajadores) from 3 to 2.
$ws.Range("C13").Value = 2

# Update worker data: period moves from 2507 to 2508 for remaining workers,
# and their "Valor Mora" / "Salario Basico" values are refreshed.

# Row 16 is now YESICA PATRICIA LOPEZ ACEVEDO (was row 17 before the delete)
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 400000
$ws.Range("G16").Value = 10000000

# Row 17 is now DALYS DANITH ESPAÑA DE ORO (was row 18 before the delete)
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500
